$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1349599
$ws.Range("C4").Value = 2290
$ws.Range("E4").Value = 1031417
$ws.Range("G4").Value = 64
$ws.Range("H4").Value = 80101

# Row 19 - Paises Bajos
$ws.Range("F19").Value = 507

# Row 59 - Kazajistan
$ws.Range("B59").Value = 5076
$ws.Range("C59").Value = 101
$ws.Range("D59").Value = 1901
$ws.Range("E59").Value = 3144

# Row 79 - Bulgaria
$ws.Range("B79").Value = 1965
$ws.Range("C79").Value = 44
$ws.Range("E79").Value = 1430
$ws.Range("F79").Value = 58
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 91

# Row 135 - Sierra Leona
$ws.Range("B135").Value = 307
$ws.Range("C135").Value = 16
$ws.Range("D135").Value = 67
$ws.Range("E135").Value = 222

# Row 147 - Birmania
$ws.Range("B147").Value = 180
$ws.Range("C147").Value = 2
$ws.Range("D147").Value = 72
$ws.Range("E147").Value = 102

# Row 166 - Mozambique
$ws.Range("B166").Value = 91
$ws.Range("C166").Value = 4
$ws.Range("E166").Value = 57
